$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E3").Value = '[Tao%Chen%NULL%0,           Di%Wu%NULL%0,           Huilong%Chen%NULL%0,           Weiming%Yan%NULL%0,           Danlei%Yang%NULL%0,           Guang%Chen%NULL%0,           Ke%Ma%NULL%0,           Dong%Xu%NULL%0,           Haijing%Yu%NULL%0,           Hongwu%Wang%NULL%0,           Tao%Wang%NULL%0,           Wei%Guo%NULL%0,           Jia%Chen%NULL%0,           Chen%Ding%NULL%0,           Xiaoping%Zhang%NULL%0,           Jiaquan%Huang%NULL%0,           Meifang%Han%NULL%0,           Shusheng%Li%NULL%0,           Xiaoping%Luo%NULL%0,           Jianping%Zhao%NULL%0,           Qin%Ning%NULL%0]'
$ws.Range("E4").Value = '[Yuan%Yu%NULL%0,           Dan%Xu%NULL%0,           Shouzhi%Fu%NULL%0,           Jun%Zhang%NULL%0,           Xiaobo%Yang%NULL%0,           Liang%Xu%NULL%0,           Jiqian%Xu%NULL%0,           Yongran%Wu%NULL%0,           Chaolin%Huang%NULL%0,           Yaqi%Ouyang%NULL%0,           Luyu%Yang%NULL%0,           Minghao%Fang%NULL%0,           Hongwen%Xiao%NULL%0,           Jing%Ma%NULL%0,           Wei%Zhu%NULL%0,           Song%Hu%NULL%0,           Quan%Hu%NULL%0,           Daoyin%Ding%NULL%0,           Ming%Hu%NULL%0,           Guochao%Zhu%NULL%0,           Weijiang%Xu%NULL%0,           Jun%Guo%NULL%0,           Jinglong%Xu%NULL%0,           Haitao%Yuan%NULL%0,           Bin%Zhang%NULL%0,           Zhui%Yu%yuzhui@whu.edu.cn%0,           Dechang%Chen%icudechangchen@163.com%0,           Shiying%Yuan%yuan_shiying@163.com%0,           You%Shang%you_shanghust@163.com%0]'
$ws.Range("E5").Value = '[Barnaby%D.P.%coreGivesNoEmail%0,          Becker%L.B.%coreGivesNoEmail%0,          Chelico%J.D.%coreGivesNoEmail%0,          Cohen%S.L.%coreGivesNoEmail%0,          Cookingham%J.%coreGivesNoEmail%0,          Coppa%K.%coreGivesNoEmail%0,          Crawford%J.M.%coreGivesNoEmail%0,          Davidson%K.W.%coreGivesNoEmail%0,          Diefenbach%M.A.%coreGivesNoEmail%0,          Dominello%A.J.%coreGivesNoEmail%0,          Duer-Hefele%J.%coreGivesNoEmail%0,          Falzon%L.%coreGivesNoEmail%0,          Gitlin%J.%coreGivesNoEmail%0,          Hajizadeh%N.%coreGivesNoEmail%0,          Harvin%T.G.%coreGivesNoEmail%0,          Hirsch%J.S.%coreGivesNoEmail%0,          Hirschwerk%D.A.%coreGivesNoEmail%0,          Kim%E.J.%coreGivesNoEmail%0,          Kozel%Z.M.%coreGivesNoEmail%0,          Marrast%L.M.%coreGivesNoEmail%0,          McGinn%T.%coreGivesNoEmail%0,          Mogavero%J.N.%coreGivesNoEmail%0,          Narasimhan%M.%coreGivesNoEmail%0,          Osorio%G.A.%coreGivesNoEmail%0,          Qiu%M.%coreGivesNoEmail%0,          Richardson%S.%coreGivesNoEmail%0,          Zanos%T.P.%coreGivesNoEmail%0]'
$ws.Range("E6").Value = '[Antonelli%Massimo%coreGivesNoEmail%0,          Cabrini%Luca%coreGivesNoEmail%0,          Castelli%Antonio%coreGivesNoEmail%0,          Cecconi%Maurizio%coreGivesNoEmail%0,          Cereda%Danilo%coreGivesNoEmail%0,          Coluccello%Antonio%coreGivesNoEmail%0,          Foti%Giuseppe%coreGivesNoEmail%0,          Fumagalli%Roberto%coreGivesNoEmail%0,          Grasselli%Giacomo%coreGivesNoEmail%0,          Iotti%Giorgio%coreGivesNoEmail%0,          Latronico%Nicola%coreGivesNoEmail%0,          Lorini%Luca%coreGivesNoEmail%0,          Merler%Stefano%coreGivesNoEmail%0,          Natalini%Giuseppe%coreGivesNoEmail%0,          Pesenti%Antonio%coreGivesNoEmail%0,          Piatti%Alessandra%coreGivesNoEmail%0,          Ranieri%Marco Vito%coreGivesNoEmail%0,          Scandroglio%Anna Mara%coreGivesNoEmail%0,          Storti%Enrico%coreGivesNoEmail%0,          Zanella%Alberto%coreGivesNoEmail%0,          Zangrillo%Alberto%coreGivesNoEmail%0]'
$ws.Range("E7").Value = '[Pavan K.%Bhatraju%NULL%0,           Bijan J.%Ghassemieh%NULL%0,           Michelle%Nichols%NULL%0,           Richard%Kim%NULL%0,           Keith R.%Jerome%NULL%0,           Arun K.%Nalla%NULL%0,           Alexander L.%Greninger%NULL%0,           Sudhakar%Pipavath%NULL%0,           Mark M.%Wurfel%NULL%0,           Laura%Evans%NULL%0,           Patricia A.%Kritek%NULL%0,           T. Eoin%West%NULL%0,           Andrew%Luks%NULL%0,           Anthony%Gerbino%NULL%0,           Chris R.%Dale%NULL%0,           Jason D.%Goldman%NULL%0,           Shane%O’Mahony%NULL%0,           Carmen%Mikacenic%NULL%0]'
$ws.Range("E8").Value = '[Hideaki%Kato%NULL%1,           Hiroyuki%Shimizu%NULL%1,           Yasushi%Shibue%NULL%1,           Tomohiro%Hosoda%NULL%1,           Keisuke%Iwabuchi%NULL%1,           Kotaro%Nagamine%NULL%1,           Hiroki%Saito%NULL%1,           Reimin%Sawada%NULL%1,           Takayuki%Oishi%NULL%1,           Jun%Tsukiji%NULL%1,           Hiroyuki%Fujita%NULL%1,           Ryosuke%Furuya%NULL%1,           Makoto%Masuda%NULL%1,           Osamu%Akasaka%NULL%1,           Yu%Ikeda%NULL%1,           Mitsuo%Sakamoto%NULL%1,           Kazuya%Sakai%NULL%1,           Munehito%Uchiyama%NULL%1,           Hiroki%Watanabe%NULL%1,           Nobuhiro%Yamaguchi%NULL%1,           Ryoko%Higa%NULL%1,           Akiko%Sasaki%NULL%1,           Katsuaki%Tanaka%NULL%1,           Yukitoshi%Toyoda%NULL%1,           Shinsuke%Hamanaka%NULL%1,           Naoki%Miyazawa%NULL%1,           Atsuko%Shimizu%NULL%1,           Fumie%Fukase%NULL%1,           Shunsuke%Iwai%NULL%1,           Yuko%Komase%NULL%1,           Tsutomu%Kawasaki%NULL%1,           Isao%Nagata%NULL%1,           Yusuke%Nakayama%NULL%1,           Tetsuhiro%Takei%NULL%1,           Katsuo%Kimura%NULL%1,           Reiko%Kunisaki%NULL%1,           Makoto%Kudo%NULL%1,           Ichiro%Takeuchi%NULL%0,           Hideaki%Nakajima%NULL%1]'
$ws.Range("E9").Value = '[Chaomin%Wu%xref no email%0,    Xiaoyan%Chen%xref no email%0,    Yanping%Cai%xref no email%0,    Jia\u2019an%Xia%xref no email%0,    Xing%Zhou%xref no email%0,    Sha%Xu%xref no email%0,    Hanping%Huang%xref no email%0,    Li%Zhang%xref no email%0,    Xia%Zhou%xref no email%0,    Chunling%Du%xref no email%0,    Yuye%Zhang%xref no email%0,    Juan%Song%xref no email%0,    Sijiao%Wang%xref no email%0,    Yencheng%Chao%xref no email%0,    Zeyong%Yang%xref no email%0,    Jie%Xu%xref no email%0,    Xin%Zhou%xref no email%0,    Dechang%Chen%xref no email%0,    Weining%Xiong%xref no email%0,    Lei%Xu%xref no email%0,    Feng%Zhou%xref no email%0,    Jinjun%Jiang%xref no email%0,    Chunxue%Bai%xref no email%0,    Junhua%Zheng%xref no email%0,    Yuanlin%Song%xref no email%0]'
$ws.Range("E10").Value = '[Xiaobo%Yang%NULL%0,           Yuan%Yu%NULL%0,           Jiqian%Xu%NULL%0,           Huaqing%Shu%NULL%0,           Jia''an%Xia%NULL%0,           Hong%Liu%NULL%0,           Yongran%Wu%NULL%0,           Lu%Zhang%NULL%0,           Zhui%Yu%NULL%0,           Minghao%Fang%NULL%0,           Ting%Yu%NULL%0,           Yaxin%Wang%NULL%0,           Shangwen%Pan%NULL%0,           Xiaojing%Zou%NULL%0,           Shiying%Yuan%NULL%0,           You%Shang%NULL%0]'
$ws.Range("E11").Value = '[Peng Peng%Xu%NULL%1,           Rong Hua%Tian%NULL%1,           Song%Luo%NULL%1,           Zi Yue%Zu%NULL%1,           Bin%Fan%NULL%1,           Xi Ming%Wang%NULL%1,           Kai%Xu%NULL%1,           Jiang Tao%Wang%NULL%1,           Juan%Zhu%NULL%1,           Ji Chan%Shi%NULL%1,           Feng%Chen%NULL%1,           Bing%Wan%NULL%1,           Zhi Han%Yan%NULL%1,           Rong Pin%Wang%NULL%1,           Wen%Chen%NULL%1,           Wen Hui%Fan%NULL%1,           Can%Zhang%NULL%1,           Meng Jie%Lu%NULL%1,           Zhi Yuan%Sun%NULL%1,           Chang Sheng%Zhou%NULL%1,           Li Na%Zhang%NULL%1,           Fei%Xia%NULL%1,           Li%Qi%NULL%1,           Wei%Zhang%NULL%0,           Jing%Zhong%NULL%1,           Xiao Xue%Liu%NULL%1,           Qi Rui%Zhang%NULL%1,           Guang Ming%Lu%NULL%1,           Long Jiang%Zhang%NULL%1]'
$ws.Range("E12").Value = '[Matthew J%Cummings%NULL%2,           Matthew R%Baldwin%NULL%0,           Darryl%Abrams%NULL%1,           Samuel D%Jacobson%NULL%1,           Benjamin J%Meyer%NULL%1,           Elizabeth M%Balough%NULL%1,           Justin G%Aaron%NULL%1,           Jan%Claassen%NULL%1,           LeRoy E%Rabbani%NULL%1,           Jonathan%Hastie%NULL%1,           Beth R%Hochman%NULL%1,           John%Salazar-Schicchi%NULL%1,           Natalie H%Yip%NULL%1,           Daniel%Brodie%NULL%2,           Max R%O''Donnell%NULL%1]'
$ws.Range("E13").Value = '[Michael G%Argenziano%NULL%0,           Samuel L%Bruce%NULL%0,           Cody L%Slater%NULL%0,           Jonathan R%Tiao%NULL%0,           Matthew R%Baldwin%NULL%0,           R Graham%Barr%NULL%0,           Bernard P%Chang%NULL%0,           Katherine H%Chau%NULL%0,           Justin J%Choi%NULL%0,           Nicholas%Gavin%NULL%0,           Parag%Goyal%NULL%0,           Angela M%Mills%NULL%0,           Ashmi A%Patel%NULL%0,           Marie-Laure S%Romney%NULL%0,           Monika M%Safford%NULL%0,           Neil W%Schluger%NULL%0,           Soumitra%Sengupta%NULL%0,           Magdalena E%Sobieszczyk%NULL%0,           Jason E%Zucker%NULL%0,           Paul A%Asadourian%NULL%0,           Fletcher M%Bell%NULL%0,           Rebekah%Boyd%NULL%0,           Matthew F%Cohen%NULL%0,           MacAlistair I%Colquhoun%NULL%0,           Lucy A%Colville%NULL%0,           Joseph H%de Jonge%NULL%0,           Lyle B%Dershowitz%NULL%0,           Shirin A%Dey%NULL%0,           Katherine A%Eiseman%NULL%0,           Zachary P%Girvin%NULL%0,           Daniella T%Goni%NULL%0,           Amro A%Harb%NULL%0,           Nicholas%Herzik%NULL%0,           Sarah%Householder%NULL%0,           Lara E%Karaaslan%NULL%0,           Heather%Lee%NULL%0,           Evan%Lieberman%NULL%0,           Andrew%Ling%NULL%0,           Ree%Lu%NULL%0,           Arthur Y%Shou%NULL%0,           Alexander C%Sisti%NULL%0,           Zachary E%Snow%NULL%0,           Colin P%Sperring%NULL%0,           Yuqing%Xiong%NULL%0,           Henry W%Zhou%NULL%0,           Karthik%Natarajan%NULL%0,           George%Hripcsak%NULL%0,           Ruijun%Chen%NULL%0]'
$ws.Range("E14").Value = '[ Anish R.%Mitra%null%0,            Nicholas A.%Fergusson%null%0,            Elisa%Lloyd-Smith%null%0,            Andrew%Wormsbecker%null%0,            Denise%Foster%null%0,            Andrei%Karpov%null%0,            Sarah%Crowe%null%0,            Greg%Haljan%null%0,            Dean R.%Chittock%null%0,            Hussein D.%Kanji%null%0,            Mypinder S.%Sekhon%null%0,            Donald E.G.%Griesdale%null%0,         Anish R.%Mitra%null%1,         Nicholas A.%Fergusson%null%1,         Elisa%Lloyd-Smith%null%1,         Andrew%Wormsbecker%null%1,         Denise%Foster%null%1,         Andrei%Karpov%null%1,         Sarah%Crowe%null%1,         Greg%Haljan%null%1,         Dean R.%Chittock%null%1,         Hussein D.%Kanji%null%1,         Mypinder S.%Sekhon%null%1,         Donald E.G.%Griesdale%null%1]'
$ws.Range("E15").Value = '[Sara C.%Auld%NULL%0,           Mark%Caridi-Scheible%NULL%0,           James M.%Blum%NULL%0,           Chad%Robichaux%NULL%0,           Colleen%Kraft%NULL%0,           Jesse T.%Jacob%NULL%0,           Craig S.%Jabaley%NULL%0,           David%Carpenter%NULL%0,           Roberta%Kaplow%NULL%0,           Alfonso C.%Hernandez-Romieu%NULL%0,           Max W.%Adelman%NULL%0,           Greg S.%Martin%NULL%0,           Craig M.%Coopersmith%NULL%0,           David J.%Murphy%NULL%0,           NULL%NULL%NULL%0,           NULL%NULL%NULL%0]'
$ws.Range("E16").Value = '[Kevin%Hur%NULL%1,           Caroline P. E.%Price%NULL%1,           Elizabeth L.%Gray%NULL%1,           Reeti K.%Gulati%NULL%1,           Matthew%Maksimoski%NULL%1,           Samuel D.%Racette%NULL%1,           Alexander L.%Schneider%NULL%1,           Ashoke R.%Khanwalkar%NULL%1]'
